$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 06:50"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("E4").Value = 177273
$ws.Range("H4").Value = 4054

# Row 22: Australia -> Australia
$ws.Range("B22").Value = 4828
$ws.Range("C22").Value = 65
$ws.Range("E22").Value = 4463

# Row 37: Pakistan -> Pakistan
$ws.Range("B37").Value = 2031
$ws.Range("C37").Value = 93
$ws.Range("E37").Value = 1929

# Row 42: India -> India
$ws.Range("D42").Value = 124
$ws.Range("E42").Value = 1238

# Row 93: Camerun -> Afganistan
$ws.Range("A93").Value = "Afganistan"
$ws.Range("B93").Value = 196
$ws.Range("C93").Value = 22
$ws.Range("E93").Value = 187
$ws.Range("H93").Value = 4

# Row 94: Oman -> Camerun
$ws.Range("A94").Value = "Camerun"
$ws.Range("B94").Value = 193
$ws.Range("D94").Value = 5
$ws.Range("E94").Value = 182
$ws.Range("F94").Value = 0
$ws.Range("H94").Value = 6

# Row 95: Cuba -> Oman
$ws.Range("A95").Value = "Oman"
$ws.Range("B95").Value = 192
$ws.Range("D95").Value = 34
$ws.Range("E95").Value = 157
$ws.Range("H95").Value = 1

# Row 96: Costa de Marfil -> Cuba
$ws.Range("A96").Value = "Cuba"
$ws.Range("B96").Value = 186
$ws.Range("D96").Value = 8
$ws.Range("E96").Value = 172
$ws.Range("F96").Value = 3
$ws.Range("H96").Value = 6

# Row 97: Senegal -> Costa de Marfil
$ws.Range("A97").Value = "Costa de Marfil"
$ws.Range("B97").Value = 179
$ws.Range("D97").Value = 7
$ws.Range("E97").Value = 171
$ws.Range("H97").Value = 1

# Row 98: Afganistan -> Senegal
$ws.Range("A98").Value = "Senegal"
$ws.Range("B98").Value = 175
$ws.Range("D98").Value = 40
$ws.Range("E98").Value = 135
$ws.Range("H98").Value = 0

# Row 142: Guam -> El Salvador
$ws.Range("A142").Value = "El Salvador"
$ws.Range("F142").Value = 5

# Row 143: El Salvador -> Guam
$ws.Range("A143").Value = "Guam"
$ws.Range("F143").Value = 0

# Row 156: Guinea Ecuatorial -> Birmania
$ws.Range("A156").Value = "Birmania"
$ws.Range("D156").Value = 0
$ws.Range("H156").Value = 1

# Row 157: Birmania -> Guinea Ecuatorial
$ws.Range("A157").Value = "Guinea Ecuatorial"
$ws.Range("D157").Value = 1
$ws.Range("H157").Value = 0

# Row 168: Surinam -> Seychelles
$ws.Range("A168").Value = "Seychelles"
$ws.Range("C168").Value = 0

# Row 169: Seychelles -> Surinam
$ws.Range("A169").Value = "Surinam"
$ws.Range("C169").Value = 1

# Row 173: Laos -> Suazilandia
$ws.Range("A173").Value = "Suazilandia"

# Row 174: Granada -> Laos
$ws.Range("A174").Value = "Laos"

# Row 175: Suazilandia -> Granada
$ws.Range("A175").Value = "Granada"

# Row 178: Mozambique -> San Cristobal y Nieves
$ws.Range("A178").Value = "San Cristobal y Nieves"

# Row 179: San Cristobal y Nieves -> Mozambique
$ws.Range("A179").Value = "Mozambique"

# Row 185: San Martin (Parte Holandesa) -> Santa Sede
$ws.Range("A185").Value = "Santa Sede"

# Row 186: Santa Sede -> San Martin (Parte Holandesa)
$ws.Range("A186").Value = "San Martin (Parte Holandesa)"

# Row 190: Montserrat -> Islas Turcas y Caicos
$ws.Range("A190").Value = "Islas Turcas y Caicos"

# Row 191: Islas Turcas y Caicos -> Montserrat
$ws.Range("A191").Value = "Montserrat"

# Row 194: Nicaragua -> Nepal
$ws.Range("A194").Value = "Nepal"
$ws.Range("D194").Value = 1
$ws.Range("H194").Value = 0

# Row 195: Nepal -> Nicaragua
$ws.Range("A195").Value = "Nicaragua"
$ws.Range("D195").Value = 0
$ws.Range("H195").Value = 1

# Row 197: Gambia -> Botsuana
$ws.Range("A197").Value = "Botsuana"

# Row 198: Botsuana -> Gambia
$ws.Range("A198").Value = "Gambia"

# Row 199: Republica de Africa Central -> Islas Virgenes Britanicas
$ws.Range("A199").Value = "Islas Virgenes Britanicas"

# Row 201: Islas Virgenes Britanicas -> Republica de Africa Central
$ws.Range("A201").Value = "Republica de Africa Central"
